$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Add new header in K1 and values in K2:K4
$ws.Range("K1").Value = "testbench"
$ws.Range("K2").Value = "SYS-110.tbc"
$ws.Range("K3").Value = "SYS-110.tbc"
$ws.Range("K4").Value = "SYS-110.tbc"

# Update selection to match the recorded state after the edit
$ws.Range("K13").Select()
